# Apply the "456a3b4" gh-pages data refresh to 苏州-漫展信息.xlsx
#
# Summary of the change:
#  - Sheet "展览"   (exhibitions):  "想去人数" (F column) counters bumped for many rows.
#  - Sheet "演出"   (performances): a brand-new event ("小鸳鸯上清联欢会") was scraped in
#                                   ahead of the existing rows, so a row is inserted at
#                                   row 2 and the two pre-existing rows shift down one.
#  - Sheet "本地生活" (local life):  untouched.
#  - Sheet "全部类型" (all types):   same F-column counter bumps as 展览/演出, plus the
#                                   row that already held the 2024-08-03 slot gets its
#                                   event info overwritten in place by the new event
#                                   (this sheet does not shift rows - it just carries
#                                   whatever event currently occupies each date slot).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "展览" - bump the "想去人数" counters.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    "F2"  = 233
    "F3"  = 1422
    "F4"  = 19993
    "F5"  = 796
    "F6"  = 310
    "F8"  = 14
    "F9"  = 7545
    "F10" = 510
    "F11" = 734
    "F12" = 260
    "F14" = 156
    "F15" = 117
    "F16" = 9
    "F17" = 234
    "F20" = 411
    "F21" = 72
    "F22" = 680
    "F23" = 49
    "F24" = 65
    "F25" = 68
    "F26" = 319
    "F27" = 1099
    "F29" = 18
    "F30" = 180
    "F31" = 5219
    "F33" = 63
    "F34" = 2832
    "F37" = 52
    "F38" = 12617
    "F39" = 1331
    "F40" = 77
    "F41" = 24
    "F43" = 260
    "F44" = 364
    "F45" = 3996
    "F46" = 319
}
foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

# ---------------------------------------------------------------------------
# 2) Sheet "演出" - insert the new 2024-08-03 event as row 2, push the rest down.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Insert a blank row above the current row 2 (shifts old rows 2,3 -> 3,4).
$ws2.Rows.Item(2).Insert()

# The freshly inserted row has no formatting; clone it from the row below
# (the old row 2, now row 3) so the index cell keeps the bold/bordered style.
$ws2.Range("A3").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

# Row index numbers (column A) for all three data rows.
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3

# New row 2 content. Force text formatting on the date cell first so Excel
# doesn't auto-convert the "yyyy-mm-dd" literal into a real date, then strip
# the formatting back off so the cell matches its plain, unstyled siblings.
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "2024-08-03"
$ws2.Range("B2").ClearFormats()

$ws2.Range("C2").Value = "苏州·小鸳鸯上清联欢会——「遇见平江」配音演员专场见面会"
$ws2.Range("D2").Value = "东苑路1号公共文化中心内 苏州保利大剧院-小剧场"
$ws2.Range("E2").Value = "2024.08.03 11:00-08.03 18:00"
$ws2.Range("F2").Value = 151
$ws2.Range("G2").Value = 398
$ws2.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=89155"
$ws2.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202407/OEGdt4u11720690236254.png"

# The event that used to be row 3 ("苏州·爱乐之城") also picked up a refreshed
# "想去人数" count while it shifted down to row 4.
$ws2.Range("F4").Value = 34

# ---------------------------------------------------------------------------
# 3) Sheet "全部类型" - same counter bumps, plus overwrite row 31 in place with
#    the new event (this combined sheet doesn't shift rows).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    "F2"  = 233
    "F3"  = 1422
    "F4"  = 19993
    "F5"  = 796
    "F6"  = 310
    "F8"  = 14
    "F9"  = 7545
    "F10" = 510
    "F11" = 734
    "F12" = 260
    "F14" = 156
    "F15" = 117
    "F16" = 9
    "F17" = 234
    "F20" = 411
    "F21" = 72
    "F22" = 680
    "F23" = 49
    "F24" = 65
    "F25" = 68
    "F26" = 319
    "F27" = 1099
    "F29" = 18
    "F30" = 180
    "F34" = 63
    "F35" = 34
    "F36" = 2832
    "F37" = 0
    "F39" = 52
    "F40" = 12617
    "F41" = 1331
    "F42" = 77
    "F43" = 24
    "F45" = 260
    "F46" = 364
    "F47" = 3996
    "F48" = 319
}
foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}

# Row 31 ("苏州·星部落动漫嘉年华", also on 2024-08-03) is replaced by the new event.
$ws4.Range("C31").Value = "苏州·小鸳鸯上清联欢会——「遇见平江」配音演员专场见面会"
$ws4.Range("D31").Value = "东苑路1号公共文化中心内 苏州保利大剧院-小剧场"
$ws4.Range("E31").Value = "2024.08.03 11:00-08.03 18:00"
$ws4.Range("F31").Value = 151
$ws4.Range("G31").Value = 398
$ws4.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=89155"
$ws4.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202407/OEGdt4u11720690236254.png"

Write-Host "Applied 456a3b4 refresh."
